$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 1359
    3 = 1977
    4 = 227
    6 = 6347
    7 = 218
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
